# Atualizacao dos dados bibi e add
# Insere um novo dia (06/2025) no topo dos dados e corrige alguns valores
# de total_venda do periodo 05/2025 que foram recalculados.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova linha logo apos o cabecalho (linha 1), empurrando todos
# os dados existentes uma linha para baixo.
$ws.Rows.Item(2).Insert()

# Remove a formatacao herdada da linha do cabecalho para que a nova linha
# de dados fique igual as demais linhas de dados (sem estilo proprio).
$ws.Range("A2:E2").ClearFormats()

# Preenche a nova linha com o novo registro de 06/2025.
$ws.Cells.Item(2, 1).Value2 = 2
$ws.Cells.Item(2, 2).Value2 = 3672.65
$ws.Cells.Item(2, 3).Value2 = 6
$ws.Cells.Item(2, 4).Value2 = 2025
$ws.Cells.Item(2, 5).Value2 = "06/2025"

# Corrige os valores de total_venda (coluna B) de alguns dias de 05/2025
# que foram atualizados (as demais colunas desses registros nao mudam).
$ws.Cells.Item(11, 2).Value2 = 34513.04
$ws.Cells.Item(15, 2).Value2 = 21165.12
$ws.Cells.Item(16, 2).Value2 = 8933.42
$ws.Cells.Item(17, 2).Value2 = 27740.2
$ws.Cells.Item(19, 2).Value2 = 27841.33
$ws.Cells.Item(20, 2).Value2 = 16252.73
$ws.Cells.Item(21, 2).Value2 = 24171.51
$ws.Cells.Item(22, 2).Value2 = 26203.72
$ws.Cells.Item(23, 2).Value2 = 18963.96
